$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# 1) Add the new "BloodGlucoseOnProfile" suite block (rows 26-28).
#    This is done before editing B20:B24 below, because row 28 is
#    cloned from row 23 while its Runmode value is still "Y".
# ------------------------------------------------------------------

# Section title row (clone row 18's layout, then replace its text)
$ws.Range("A18").Copy($ws.Range("A26"))
$ws.Range("A26").Value = "BloodGlucoseOnProfile"

# Column header row (clone row 19's header cells A:H)
$ws.Range("A19:H19").Copy($ws.Range("A27"))

# Data row (clone row 23's layout/values, then override the test id and description)
$ws.Range("A23:H23").Copy($ws.Range("A28"))
$ws.Range("A28").Value = "TC09"
$ws.Range("C28").Value = "Verify that, system updates BG values on Profile page"

# Hyperlinks for the new data row (UserName / Password columns)
$ws.Hyperlinks.Add($ws.Range("D28"), "mailto:neil@peter.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E28"), "mailto:Test@123") | Out-Null

# Re-apply the formatting from the row-20 equivalent cells so the
# hyperlink insertion does not leave a stray/duplicated style behind.
$ws.Range("D20:E20").Copy()
$ws.Range("D28").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# ------------------------------------------------------------------
# 2) Flip the Runmode column for the BloodGlucoseTest suite (rows 20-24)
#    from "Y" to "N".
# ------------------------------------------------------------------
$ws.Range("B20").Value = "N"
$ws.Range("B21").Value = "N"
$ws.Range("B22").Value = "N"
$ws.Range("B23").Value = "N"
$ws.Range("B24").Value = "N"

# ------------------------------------------------------------------
# 3) Restore the selection/active cell recorded in the workbook.
# ------------------------------------------------------------------
$ws.Range("B20:B24").Select() | Out-Null

Write-Host "done"
